# "testing V1G and V2G" - add a UC_ELC-Nuc_BuildRate constraint row (with a
# TFM_INS helper table) to the UC_NcapBuildRate sheet, fix the ~UC_T tag
# name, retarget the ~UC_Sets line to regions (R_S) instead of time (R_E),
# and swap the Wind/Solar build-rate figures between the two existing rows.

$wb = $excel.ActiveWorkbook

$wsBuild  = $wb.Worksheets.Item("UC_NcapBuildRate")
$wsGrowth = $wb.Worksheets.Item("UC_Growth")

# --- UC_NcapBuildRate: prep row 8 by copying row 7's formatting down -----
# (pure formatting copy - values get overwritten below, so it doesn't
# disturb the shared-string table)
$wsBuild.Range("B7:L7").Copy()
$wsBuild.Range("B8:L8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- header/tag fixups -----------------------------------------------------
$wsBuild.Range("B2").Value = "~UC_Sets: R_S: AllRegions"
$wsBuild.Range("B3").Value = ""
$wsBuild.Range("D4").Value = "~UC_T: UC_RHST~UP"

# --- new row 8 - Nuclear build rate UC -------------------------------------
$wsBuild.Range("B8").Value = "UC_ELC-Nuc_BuildRate"
$wsBuild.Range("C8").Value = "ELCNUC"
$wsBuild.Range("D8").Value = "NCAP,BUILDUP"
$wsBuild.Range("E8").Value = 1
$wsBuild.Range("F8").Value = 1
$wsBuild.Range("G8").Value = 2
$wsBuild.Range("H8").Value = ""
$wsBuild.Range("I8").Value = ""
$wsBuild.Range("J8").Value = 3
$wsBuild.Range("K8").Value = 5
$wsBuild.Range("L8").Value = "Max Nuclear capacity installed per year"

# --- new ~TFM_INS helper table in columns O:T ------------------------------
$wsBuild.Range("O4").Value = "~TFM_INS"

$wsBuild.Range("O5").Value = "attribute"
$wsBuild.Range("P5").Value = "cset_cn"
$wsBuild.Range("Q5").Value = "pset_pn"

$wsBuild.Range("Q6").Value = "IMP*Z"
$wsBuild.Range("P6").Value = "*build*"
$wsBuild.Range("O6").Value = "FLO_BND"

$wsBuild.Range("R5").Value = "year"
$wsBuild.Range("S5").Value = "value"
$wsBuild.Range("T5").Value = "limtype"

$wsBuild.Range("R6").Value = 0
$wsBuild.Range("S6").Value = 2
$wsBuild.Range("T6").Value = "UP"

# --- swap the Wind (row6) / Solar (row7) build rates -----------------------
$wsBuild.Range("F6").Value = 1
$wsBuild.Range("G6").Value = 5
$wsBuild.Range("H6").Value = 7
$wsBuild.Range("I6").Value = 10
$wsBuild.Range("J6").Value = 20

$wsBuild.Range("F7").Value = 2
$wsBuild.Range("G7").Value = 10
$wsBuild.Range("H7").Value = 15
$wsBuild.Range("I7").Value = 20
$wsBuild.Range("J7").Value = 30

# --- selections / active sheet ---------------------------------------------
$wsGrowth.Activate()
$wsGrowth.Range("I10").Select()

$wsBuild.Activate()
$wsBuild.Range("A6").Select()
